$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.8
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.73
$ws.Range("W2").Value = 6
$ws.Range("Z2").Value = 15
$ws.Range("AE2").Value = 19
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 15
$ws.Range("AM2").Value = 41
$ws.Range("AO2").Value = 10
$ws.Range("AT2").Value = 2.63
$ws.Range("AW2").Value = 6
$ws.Range("AX2").Value = 26
$ws.Range("AY2").Value = 34
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 301
